$d = $word.ActiveDocument

# Locate the "Thank you for using Pizza Order!" paragraph, then insert the
# new "Test Data" / "Input Validation" section right before it.
$targetIdx = 0
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -like "*Thank you for using Pizza Order!*") {
        $targetIdx = $idx
        break
    }
}

$prev = $d.Paragraphs.Item($targetIdx - 1)
$r = $prev.Range
$r.Collapse(0)

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="644"/>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">Test Data:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:highlight w:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">Cheese Pizza</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> – Automatically Populates Choices for Cheese Pizza</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:highlight w:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">Pepperoni Pizza</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> – Automatically Populates Choices for Pepperoni Pizza</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
      </w:r>
      <w:r/>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">Results:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">Item: Cheese Pizza</w:t>
      </w:r>
      <w:r/>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">Crust Choice: Hand Tossed - $1</w:t>
      </w:r>
      <w:r/>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">Sauce Choice: Tasty Classic - $0.5</w:t>
      </w:r>
      <w:r/>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">Cheese Choice: Extra Cheese - $1.5</w:t>
      </w:r>
      <w:r/>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">Size Choice: Medium Pizza 12 Inch - $4</w:t>
      </w:r>
      <w:r/>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">Topping 1 Choice: None - $0</w:t>
      </w:r>
      <w:r/>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">Topping 2 Choice: None - $0</w:t>
      </w:r>
      <w:r/>
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
      </w:r>
      <w:r/>
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">Total : $7.0</w:t>
      </w:r>
      <w:r/>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
      </w:r>
      <w:r/>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">Item: Pepperoni Pizza</w:t>
      </w:r>
      <w:r/>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">Crust Choice: Hand Tossed - $1</w:t>
      </w:r>
      <w:r/>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">Sauce Choice: Bold Marinara Sauce - $0.5</w:t>
      </w:r>
      <w:r/>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">Cheese Choice: Normal Cheese - $1</w:t>
      </w:r>
      <w:r/>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">Size Choice: Medium Pizza 12 Inch - $4</w:t>
      </w:r>
      <w:r/>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">Topping 1 Choice: Pepperoni - $1</w:t>
      </w:r>
      <w:r/>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">Topping 2 Choice: None - $0</w:t>
      </w:r>
      <w:r/>
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
      </w:r>
      <w:r/>
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">Total : $7.5</w:t>
      </w:r>
      <w:r/>
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
      </w:r>
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="644"/>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">Input Validation:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">I had to validate that an item had been selected before anything could be added to cart. This application does not add an order to cart if nothing has been selected.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> Once something is selected, the options are populated automatically for the order, and can be further customized.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
      </w:r>
      <w:r>
        <w:rPr>
          <w:highlight w:val="none"/>
        </w:rPr>
      </w:r>
    </w:p>
'@
$r.InsertXML($xml)

# Grow the page height to accommodate the new content (matches the
# sectPr/pgSz change in the target revision).
$ps = $d.PageSetup
$ps.PageHeight = 2267.7
